$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86; this shifts the existing rows 86..134
# down to 87..135 (preserving all of their values/formatting).
$ws.Rows(86).Insert()

# Populate the newly inserted row 86 with the new record.
$ws.Cells.Item(86, 1).Value = 4
$ws.Cells.Item(86, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(86, 3).Value = 'Los Lagos'
$ws.Cells.Item(86, 4).Value = 44806
$ws.Cells.Item(86, 5).Value = 10
$ws.Cells.Item(86, 6).Value = 100112022
$ws.Cells.Item(86, 7).Value = 'Arveja Verde'
$ws.Cells.Item(86, 8).Value = 'Perfection'
$ws.Cells.Item(86, 9).Value = 'Primera'
$ws.Cells.Item(86, 10).Value = 70
$ws.Cells.Item(86, 11).Value = 43000
$ws.Cells.Item(86, 12).Value = 43000
$ws.Cells.Item(86, 13).Value = 43000
$ws.Cells.Item(86, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(86, 15).Value = 'Provincia de Huasco'
$ws.Cells.Item(86, 16).Value = 1720
$ws.Cells.Item(86, 17).Value = 25
$ws.Cells.Item(86, 18).Value = 'Hortaliza'
